$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.664.90'
$ws.Range("E2").Value = '  +2.52%  '

# Row 3
$ws.Range("D3").Value = '3.833.99'
$ws.Range("E3").Value = '  +1.17%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '632.38'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +5.19%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.35'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.74%  '

# Row 7
$ws.Range("D7").Value = '3.830.09'
$ws.Range("E7").Value = '  +1.15%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.18%  '

# Row 9
$ws.Range("E9").Value = '  +0.83%  '

# Row 10
$ws.Range("E10").Value = '  +2.30%  '

# Row 11
$ws.Range("E11").Value = '  +1.02%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.67'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.56%  '

# Row 13
$ws.Range("E13").Value = '  +0.86%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.13'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.19%  '

# Row 15
$ws.Range("D15").Value = '4.473.27'
$ws.Range("E15").Value = '  +1.01%  '

# Row 16
$ws.Range("D16").Value = '3.880.03'
$ws.Range("E16").Value = '  +2.99%  '

# Row 17
$ws.Range("D17").Value = '69.606.13'
$ws.Range("E17").Value = '  +2.36%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.22'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.58%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.16'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.27%  '

# Row 20
$ws.Range("E20").Value = '  -0.18%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '469.40'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.55%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.73'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.27%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.710'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.97%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000152'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.59%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.88'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.19'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.10%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.03'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.06%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.12'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.39%  '

# Row 29
$ws.Range("E29").Value = '  +0.12%  '

# Row 30
$ws.Range("D30").Value = '3.979.47'
$ws.Range("E30").Value = '  +1.02%  '

# Row 31
$ws.Range("E31").Value = '  +2.37%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.23'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.77%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.31'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.35%  '

# Row 34
$ws.Range("E34").Value = '  +0.17%  '

# Row 35
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.776.45'
$ws.Range("E35").Value = '  +0.95%  '

# Row 36
$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.19%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.07'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.03%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.104'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.17%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.150'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +7.92%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.46'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +6.35%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.93'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.68%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.983'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.53%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.17%  '

# Row 44
$ws.Range("E44").Value = '  +0.01%  '

# Row 45
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '156.13'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.56%  '

# Row 46
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.303'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +1.49%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.84'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.71%  '

# Row 48
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.42'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +5.65%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '47.03'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.84%  '

# Row 50
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.93'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.68%  '

# Row 51
$ws.Range("E51").Value = '  +1.53%  '
